# Generate Report for Handoff
#
# The localization-status report tracks, per source file, the most recent
# handoff timestamp both on the "Overview" sheet (one row per file, a single
# "Latest Handoff Date" column) and on each per-locale sheet ("zh-cn",
# "de-de" — one row per file, a "Latest Handoff Datetime" column for that
# locale's most recent handoff file).
#
# A new handoff was generated for the file
# "e6ff6c56-4918-48d4-ab0f-c5c4e63ca0c0.md", so its handoff timestamps are
# refreshed on all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for e6ff6c56-4918-48d4-ab0f-c5c4e63ca0c0.md (row 6) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-30-17 18:30:55"

# --- zh-cn sheet: row for e6ff6c56-4918-48d4-ab0f-c5c4e63ca0c0.md (row 6) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-17 18:30:51"

# --- de-de sheet: row for e6ff6c56-4918-48d4-ab0f-c5c4e63ca0c0.md (row 6) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-17 18:30:55"
